$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose B:AD (columns 2..30) contents must be swapped, keeping
# column A (the running index) untouched on each row.
$pairs = @(
    @(61, 62),
    @(88, 89),
    @(190, 191),
    @(262, 263),
    @(302, 305),
    @(303, 306)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range($ws.Cells.Item($r1, 2), $ws.Cells.Item($r1, 30))
    $range2 = $ws.Range($ws.Cells.Item($r2, 2), $ws.Cells.Item($r2, 30))

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value = $vals2
    $range2.Value = $vals1
}
